# Apply the diff: version bump + minor punctuation/wording fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 1.0 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text: fix "usuario" -> "usuário" and add trailing period
# Appears once per test case block (TC1..TC4)
$precondition = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B8").Value = $precondition
$ws.Range("B18").Value = $precondition
$ws.Range("B28").Value = $precondition
$ws.Range("B38").Value = $precondition

# Step 1 text: add trailing period
$step1 = "Chefe Acessa a funcionalidade Minha Conta Bancária (menu)."
$ws.Range("B10").Value = $step1
$ws.Range("B20").Value = $step1
$ws.Range("B30").Value = $step1
$ws.Range("B40").Value = $step1

# Step 2 expected result text: add trailing period
$step2Result = "SYSTEM Apresenta os campos (banco/agência/conta corrente) alterados."
$ws.Range("D11").Value = $step2Result
$ws.Range("D21").Value = $step2Result
$ws.Range("D31").Value = $step2Result

# TC4 step 2 expected result: "conta bancários" -> "conta bancária"
$ws.Range("D41").Value = "SYSTEM Exibe mensagens informativas (MSG403 - Informativos sobre a atualização de conta bancária (dados bancários)) para o usuário sobre a manutenção de informações bancárias."
